# Auto-generated edit script: updates currentAveragePrice-derived profit columns
# (H..N) for specific rows across multiple sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 547.1667
$ws.Range("I28").Value = 601.1
$ws.Range("K28").Value = 601.1
$ws.Range("M28").Value = -116.1
$ws.Range("H33").Value = 13998.556
$ws.Range("I33").Value = 16771.334
$ws.Range("K33").Value = 16771.334
$ws.Range("M33").Value = -16542.334
$ws.Range("H106").Value = 3964.7273
$ws.Range("I106").Value = 4116.3
$ws.Range("J106").Value = 2449
$ws.Range("K106").Value = 4116.3
$ws.Range("L106").Value = 2449
$ws.Range("M106").Value = -3485.3
$ws.Range("N106").Value = -3711
$ws.Range("H132").Value = 9767.034
$ws.Range("I132").Value = 12056.546
$ws.Range("K132").Value = 36169.638
$ws.Range("M132").Value = -33639.638

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5270228
$ws.Range("I61").Value = 8909.362999999999
$ws.Range("J61").Value = 12504541
$ws.Range("K61").Value = 8909.362999999999
$ws.Range("L61").Value = 12504541
$ws.Range("M61").Value = -8697.362999999999
$ws.Range("N61").Value = -12504965
$ws.Range("H122").Value = 20834880
$ws.Range("J122").Value = 2005.25
$ws.Range("L122").Value = 6015.75
$ws.Range("N122").Value = -10915.75
$ws.Range("H132").Value = 405775.2
$ws.Range("I132").Value = 474201.28
$ws.Range("K132").Value = 1422603.84
$ws.Range("M132").Value = -1420073.84
$ws.Range("H136").Value = 5270228
$ws.Range("I136").Value = 8909.362999999999
$ws.Range("J136").Value = 12504541
$ws.Range("K136").Value = 26728.089
$ws.Range("L136").Value = 37513623
$ws.Range("M136").Value = -24178.089
$ws.Range("N136").Value = -37518723

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6166.2856
$ws.Range("I105").Value = 5368.933
$ws.Range("J105").Value = 8159.6665
$ws.Range("K105").Value = 5368.933
$ws.Range("L105").Value = 8159.6665
$ws.Range("M105").Value = -3621.933
$ws.Range("N105").Value = -11653.6665
$ws.Range("H134").Value = 3092674.8
$ws.Range("I134").Value = 4530.3335
$ws.Range("K134").Value = 13591.0005
$ws.Range("M134").Value = -11056.0005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 88413.78
$ws.Range("I16").Value = 72847.78999999999
$ws.Range("J16").Value = 112627.555
$ws.Range("K16").Value = 72847.78999999999
$ws.Range("L16").Value = 112627.555
$ws.Range("M16").Value = -72560.78999999999
$ws.Range("N16").Value = -113201.555
$ws.Range("H22").Value = 795.6445
$ws.Range("I22").Value = 373.16217
$ws.Range("J22").Value = 2749.625
$ws.Range("K22").Value = 373.16217
$ws.Range("L22").Value = 2749.625
$ws.Range("M22").Value = -23.16217
$ws.Range("N22").Value = -3449.625
$ws.Range("H31").Value = 1324869.6
$ws.Range("J31").Value = 2222.7
$ws.Range("L31").Value = 2222.7
$ws.Range("N31").Value = -2812.7
$ws.Range("H34").Value = 1324869.6
$ws.Range("J34").Value = 2222.7
$ws.Range("L34").Value = 2222.7
$ws.Range("N34").Value = -2626.7
$ws.Range("H113").Value = 88413.78
$ws.Range("I113").Value = 72847.78999999999
$ws.Range("J113").Value = 112627.555
$ws.Range("K113").Value = 72847.78999999999
$ws.Range("L113").Value = 112627.555
$ws.Range("M113").Value = -70677.78999999999
$ws.Range("N113").Value = -116967.555

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 155.57143
$ws.Range("I7").Value = 167.8
$ws.Range("K7").Value = 503.4
$ws.Range("M7").Value = -391.4
$ws.Range("H107").Value = 340.35
$ws.Range("J107").Value = 321.4
$ws.Range("L107").Value = 964.1999999999999
$ws.Range("N107").Value = -4804.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 48965
$ws.Range("J74").Value = 48965
$ws.Range("L74").Value = 48965
$ws.Range("N74").Value = -50837
$ws.Range("H77").Value = 48965
$ws.Range("J77").Value = 48965
$ws.Range("L77").Value = 146895
$ws.Range("N77").Value = -156255
$ws.Range("H126").Value = 7891.7856
$ws.Range("I126").Value = 11434
$ws.Range("J126").Value = 3168.8333
$ws.Range("K126").Value = 34302
$ws.Range("L126").Value = 9506.499899999999
$ws.Range("M126").Value = -31832
$ws.Range("N126").Value = -14446.4999
$ws.Range("H132").Value = 11750.632
$ws.Range("I132").Value = 9812.735000000001
$ws.Range("J132").Value = 28222.75
$ws.Range("K132").Value = 29438.205
$ws.Range("L132").Value = 84668.25
$ws.Range("M132").Value = -26908.205
$ws.Range("N132").Value = -89728.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6352.2856
$ws.Range("I7").Value = 8020.7144
$ws.Range("J7").Value = 4683.857
$ws.Range("K7").Value = 8020.7144
$ws.Range("L7").Value = 4683.857
$ws.Range("M7").Value = -7908.7144
$ws.Range("N7").Value = -4907.857
$ws.Range("H46").Value = 4614.4165
$ws.Range("I46").Value = 1633.3334
$ws.Range("J46").Value = 5608.1113
$ws.Range("K46").Value = 1633.3334
$ws.Range("L46").Value = 5608.1113
$ws.Range("M46").Value = -1445.3334
$ws.Range("N46").Value = -5984.1113
$ws.Range("H69").Value = 24999.5
$ws.Range("J69").Value = 24999.5
$ws.Range("L69").Value = 24999.5
$ws.Range("N69").Value = -26621.5
$ws.Range("H72").Value = 24999.5
$ws.Range("J72").Value = 24999.5
$ws.Range("L72").Value = 74998.5
$ws.Range("N72").Value = -83110.5
$ws.Range("H100").Value = 2992
$ws.Range("I100").Value = 1986
$ws.Range("K100").Value = 1986
$ws.Range("M100").Value = -1445
$ws.Range("H126").Value = 6352.2856
$ws.Range("I126").Value = 8020.7144
$ws.Range("J126").Value = 4683.857
$ws.Range("K126").Value = 24062.1432
$ws.Range("L126").Value = 14051.571
$ws.Range("M126").Value = -21592.1432
$ws.Range("N126").Value = -18991.571
$ws.Range("H136").Value = 5958459.5
$ws.Range("I136").Value = 4811211
$ws.Range("K136").Value = 14433633
$ws.Range("M136").Value = -14431083

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 27666.666
$ws.Range("J54").Value = 27666.666
$ws.Range("L54").Value = 27666.666
$ws.Range("N54").Value = -28706.666
$ws.Range("H106").Value = 30000
$ws.Range("I106").Value = 30000
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 30000
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("M106").Value = -28738
